$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four "special parking" headers (columns J:M) were mislabeled - rename
# them to match the "Anzahl Stellplätze ..." naming used by the neighboring
# "Anzahl Stellplätze" (col I) header, fixing the missing/garbled special
# parking capacity headers.
$ws.Range("J1").Value2 = "Anzahl Stellplätze Carsharing"
$ws.Range("K1").Value2 = "Anzahl Stellplätze Lademöglichkeit"
$ws.Range("L1").Value2 = "Anzahl Stellplätze Frauen"
$ws.Range("M1").Value2 = "Anzahl Stellplätze Behinderte"

# Resize columns I and J so the new header text fits.
$ws.Columns("I").ColumnWidth = 15.7265625
$ws.Columns("J").ColumnWidth = 25.54296875

# Set up the page for printing (A4 portrait).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Leave the selection where the editor ended up after the edit.
$ws.Range("I10").Select()
